# Update Active_Outages.xlsx - 6/18/2025, 4:28:06 PM
# Updates the "Elapsed Duration(Hrs)" (column G) values across the regional
# sheets to reflect newly recalculated durations, and appends a new outage
# row (row 5, duplicating row 4) on sheet "R1".

$wb = $excel.ActiveWorkbook

# --- Sheet R1 ---
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3929:42:11"
$ws1.Range("G3").Value = "69:14:49"

# Append new row 5, duplicating row 4's content exactly.
$ws1.Range("A5").Value = ""
$ws1.Range("B5").Value = "R4"
$ws1.Range("C5").Value = ""
$ws1.Range("D5").Value = "JED0123"
$ws1.Range("E5").Value = ""
$ws1.Range("F5").Value = ""
$ws1.Range("G5").Value = ""
$ws1.Range("H5").Value = ""
$ws1.Range("I5").Value = "SCECO"
$ws1.Range("J5").Value = "In progress"
$ws1.Range("K5").Value = ""
$ws1.Range("L5").Value = "Latis"

# --- Sheet R2 ---
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12111:05:50"
$ws2.Range("G3").Value = "3240:49:19"
$ws2.Range("G4").Value = "479:00:53"

# --- Sheet R4 ---
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2956:55:39"
$ws4.Range("G3").Value = "184:07:54"
$ws4.Range("G4").Value = "72:20:19"
$ws4.Range("G5").Value = "69:57:52"

# --- Sheet R5 ---
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "430:54:38"

# --- Sheet R6 ---
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "71:26:56"
